# Regenerate the "K" column (column G) values in the save_data sheet.
# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals" -- the K (Strike count) column is recomputed
# and rewritten for every data row (rows 2-70).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column G ("K")
$newK = @{
    2  = 1
    3  = 2
    4  = 3
    5  = 3
    6  = 1
    7  = 3
    8  = 0
    9  = 0
    10 = 1
    11 = 1
    12 = 1
    13 = 4
    14 = 0
    15 = 2
    16 = 0
    17 = 2
    18 = 0
    19 = 1
    20 = 4
    21 = 3
    22 = 0
    23 = 0
    24 = 1
    25 = 1
    26 = 3
    27 = 1
    28 = 0
    29 = 1
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 0
    36 = 0
    37 = 1
    38 = 0
    39 = 0
    40 = 1
    41 = 3
    42 = 1
    43 = 0
    44 = 2
    45 = 0
    46 = 3
    47 = 2
    48 = 0
    49 = 1
    50 = 0
    51 = 2
    52 = 0
    53 = 0
    54 = 3
    55 = 1
    56 = 2
    57 = 0
    58 = 3
    59 = 2
    60 = 1
    61 = 1
    62 = 3
    63 = 0
    64 = 2
    65 = 1
    66 = 2
    67 = 0
    68 = 1
    69 = 2
    70 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
